$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5808.875
$ws.Range("I51").Value = 1250
$ws.Range("J51").Value = 6460.143
$ws.Range("K51").Value = 1250
$ws.Range("L51").Value = 6460.143
$ws.Range("M51").Value = -766
$ws.Range("N51").Value = -7428.143
$ws.Range("H53").Value = 4873.2607
$ws.Range("I53").Value = 339.06668
$ws.Range("K53").Value = 339.06668
$ws.Range("M53").Value = 297.93332
$ws.Range("H57").Value = 67638.5
$ws.Range("J57").Value = 67638.5
$ws.Range("L57").Value = 202915.5
$ws.Range("N57").Value = -203913.5
$ws.Range("H62").Value = 4467.9375
$ws.Range("I62").Value = 1561.5
$ws.Range("K62").Value = 1561.5
$ws.Range("M62").Value = -937.5
$ws.Range("H65").Value = 4467.9375
$ws.Range("I65").Value = 1561.5
$ws.Range("K65").Value = 7807.5
$ws.Range("M65").Value = -4687.5
$ws.Range("H70").Value = 4488.222
$ws.Range("I70").Value = 861.125
$ws.Range("J70").Value = 7389.9
$ws.Range("K70").Value = 2583.375
$ws.Range("L70").Value = 22169.7
$ws.Range("M70").Value = -2313.375
$ws.Range("N70").Value = -22709.7
$ws.Range("H73").Value = 4488.222
$ws.Range("I73").Value = 861.125
$ws.Range("J73").Value = 7389.9
$ws.Range("K73").Value = 2583.375
$ws.Range("L73").Value = 22169.7
$ws.Range("M73").Value = -1647.375
$ws.Range("N73").Value = -24041.7
$ws.Range("H99").Value = 2288
$ws.Range("J99").Value = 1000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 975792.0600000001
$ws.Range("I2").Value = 1088201.9
$ws.Range("K2").Value = 1088201.9
$ws.Range("M2").Value = -1088088.9
$ws.Range("H61").Value = 2607.7144
$ws.Range("I61").Value = 2499.6553
$ws.Range("K61").Value = 2499.6553
$ws.Range("M61").Value = -2287.6553
$ws.Range("H116").Value = 975792.0600000001
$ws.Range("I116").Value = 1088201.9
$ws.Range("K116").Value = 1088201.9
$ws.Range("M116").Value = -1085907.9
$ws.Range("H132").Value = 1607.1951
$ws.Range("I132").Value = 969.4722
$ws.Range("K132").Value = 2908.4166
$ws.Range("M132").Value = -378.4166
$ws.Range("H136").Value = 2607.7144
$ws.Range("I136").Value = 2499.6553
$ws.Range("K136").Value = 7498.965899999999
$ws.Range("M136").Value = -4948.965899999999
$ws.Range("H139").Value = 83000
$ws.Range("J139").Value = 83000
$ws.Range("L139").Value = 83000
$ws.Range("N139").Value = -93280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 975792.0600000001
$ws.Range("I3").Value = 1088201.9
$ws.Range("K3").Value = 1088201.9
$ws.Range("M3").Value = -1088087.9
$ws.Range("H81").Value = 13473.111
$ws.Range("J81").Value = 13473.111
$ws.Range("L81").Value = 13473.111
$ws.Range("N81").Value = -15595.111
$ws.Range("H84").Value = 13473.111
$ws.Range("J84").Value = 13473.111
$ws.Range("L84").Value = 40419.333
$ws.Range("N84").Value = -51027.333
$ws.Range("H134").Value = 3753.4055
$ws.Range("I134").Value = 878.7059
$ws.Range("K134").Value = 2636.1177
$ws.Range("M134").Value = -101.1177000000002
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H43").Value = 32500
$ws.Range("J43").Value = 32500
$ws.Range("L43").Value = 32500
$ws.Range("N43").Value = -32868
$ws.Range("H58").Value = 5846.778
$ws.Range("I58").Value = 6566.407
$ws.Range("K58").Value = 6566.407
$ws.Range("M58").Value = -6363.407
$ws.Range("H99").Value = 3733.7856
$ws.Range("I99").Value = 3893.4
$ws.Range("J99").Value = 3334.75
$ws.Range("K99").Value = 3893.4
$ws.Range("L99").Value = 3334.75
$ws.Range("M99").Value = -2395.4
$ws.Range("N99").Value = -6330.75
$ws.Range("H101").Value = 32500
$ws.Range("J101").Value = 32500
$ws.Range("L101").Value = 32500
$ws.Range("N101").Value = -38990
$ws.Range("H107").Value = 1381.4419
$ws.Range("I107").Value = 1603.8438
$ws.Range("K107").Value = 1603.8438
$ws.Range("M107").Value = 316.1561999999999
$ws.Range("H126").Value = 3733.7856
$ws.Range("I126").Value = 3893.4
$ws.Range("J126").Value = 3334.75
$ws.Range("K126").Value = 11680.2
$ws.Range("L126").Value = 10004.25
$ws.Range("M126").Value = -9210.200000000001
$ws.Range("N126").Value = -14944.25
$ws.Range("H134").Value = 31388.451
$ws.Range("I134").Value = 38971.75
$ws.Range("J134").Value = 5388.5713
$ws.Range("K134").Value = 116915.25
$ws.Range("L134").Value = 16165.7139
$ws.Range("M134").Value = -114380.25
$ws.Range("N134").Value = -21235.7139
$ws.Range("H136").Value = 5846.778
$ws.Range("I136").Value = 6566.407
$ws.Range("K136").Value = 19699.221
$ws.Range("M136").Value = -17149.221
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2570.0454
$ws.Range("I14").Value = 2570.0454
$ws.Range("K14").Value = 7710.1362
$ws.Range("M14").Value = -7537.1362
$ws.Range("H94").Value = 10000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8336903.5
$ws.Range("I70").Value = 10003601
$ws.Range("J70").Value = 3413
$ws.Range("K70").Value = 10003601
$ws.Range("L70").Value = 3413
$ws.Range("M70").Value = -10003331
$ws.Range("N70").Value = -3953
$ws.Range("H73").Value = 8336903.5
$ws.Range("I73").Value = 10003601
$ws.Range("J73").Value = 3413
$ws.Range("K73").Value = 10003601
$ws.Range("L73").Value = 3413
$ws.Range("M73").Value = -10002665
$ws.Range("N73").Value = -5285
$ws.Range("H100").Value = 37940.5
$ws.Range("J100").Value = 37940.5
$ws.Range("L100").Value = 37940.5
$ws.Range("N100").Value = -40104.5
$ws.Range("H109").Value = 55165
$ws.Range("J109").Value = 55165
$ws.Range("L109").Value = 55165
$ws.Range("N109").Value = -57245
$ws.Range("H132").Value = 2619.0789
$ws.Range("I132").Value = 2534.9375
$ws.Range("K132").Value = 7604.8125
$ws.Range("M132").Value = -5074.8125
$ws.Range("H136").Value = 11303.667
$ws.Range("J136").Value = 11303.667
$ws.Range("L136").Value = 33911.001
$ws.Range("N136").Value = -39011.001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3434.9
$ws.Range("I46").Value = 3400
$ws.Range("J46").Value = 3437.3928
$ws.Range("K46").Value = 3400
$ws.Range("L46").Value = 3437.3928
$ws.Range("M46").Value = -3212
$ws.Range("N46").Value = -3813.3928
$ws.Range("H56").Value = 8150.25
$ws.Range("I56").Value = 28000
$ws.Range("J56").Value = 1533.6666
$ws.Range("K56").Value = 28000
$ws.Range("L56").Value = 1533.6666
$ws.Range("M56").Value = -27309
$ws.Range("N56").Value = -2915.6666
$ws.Range("H100").Value = 3089.9355
$ws.Range("I100").Value = 2252.2354
$ws.Range("J100").Value = 4107.143
$ws.Range("K100").Value = 2252.2354
$ws.Range("L100").Value = 4107.143
$ws.Range("M100").Value = -1711.2354
$ws.Range("N100").Value = -5189.143
$ws.Range("H132").Value = 5944.5537
$ws.Range("I132").Value = 6034.0713
$ws.Range("J132").Value = 5387.5557
$ws.Range("K132").Value = 18102.2139
$ws.Range("L132").Value = 16162.6671
$ws.Range("M132").Value = -15572.2139
$ws.Range("N132").Value = -21222.6671
$ws.Range("H136").Value = 35587.266
$ws.Range("I136").Value = 74657
$ws.Range("K136").Value = 223971
$ws.Range("M136").Value = -221421

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7578.7393
$ws.Range("I62").Value = 3688.25
$ws.Range("J62").Value = 7949.2617
$ws.Range("K62").Value = 3688.25
$ws.Range("L62").Value = 7949.2617
$ws.Range("M62").Value = -3064.25
$ws.Range("N62").Value = -9197.261699999999
$ws.Range("H65").Value = 7578.7393
$ws.Range("I65").Value = 3688.25
$ws.Range("J65").Value = 7949.2617
$ws.Range("K65").Value = 18441.25
$ws.Range("L65").Value = 39746.3085
$ws.Range("M65").Value = -15321.25
$ws.Range("N65").Value = -45986.3085
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
$ws.Range("H132").Value = 22468328
$ws.Range("I132").Value = 25645430
$ws.Range("J132").Value = 1817154.4
$ws.Range("K132").Value = 76936290
$ws.Range("L132").Value = 5451463.199999999
$ws.Range("M132").Value = -76933760
$ws.Range("N132").Value = -5456523.199999999
